$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ANATOMY session 1 "Recorded By" list reordered
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 3 - ANATOMY session 2 "Recorded By" list reordered
$ws.Range("G3").Value = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, majorelle.magdy@med.asu.edu.eg"

# Row 4 - ANATOMY session 3 "Recorded By" list reordered
$ws.Range("G4").Value = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# Row 5 - ANATOMY session 4 "Recorded By" list reordered
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 6 - ANATOMY session 5 status changed
$ws.Range("I6").Value = "Not Recorded"

# Row 7 - BIOCHEMISTRY LAB/CBL session 1 "Recorded By" list reordered
$ws.Range("G7").Value = "AbeerRagheb@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg"

# Class Statistics: Missing Sessions / Pending Sessions updated
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 10

# Row 9 - HISTOLOGY session 1 "Recorded By" list reordered
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# Row 12 - MICROBIOLOGY session 1 "Recorded By" list reordered
$ws.Range("G12").Value = "dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"

# Row 15 - PARASITOLOGY session 2 "Recorded By" list reordered
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Group Statistics: Missing / Pending updated for Year 2 / C1
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 10

# Row 28 - PHYSIOLOGY session 1 "Recorded By" list reordered
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
